# Update column F ("dSF") values for specific rows per repulled data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    7  = 0
    9  = -1
    10 = 3
    16 = 1
    22 = -3
    31 = -5
    33 = -3
    34 = 4
    36 = 1
    40 = 1
    41 = 1
    45 = 2
    46 = -2
    51 = 2
    52 = 2
    53 = -1
    60 = -2
    61 = -1
    62 = 0
    66 = 0
    72 = 2
    74 = 2
    76 = 4
    78 = -3
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
